$wb = $excel.ActiveWorkbook

# --- Sheet "List1": update measurement timestamp for row 15 ---
# (this is the raw source timestamp; columns B/C/E on this sheet recalc
#  automatically because they hold shared formulas dependent on column A)
$ws1 = $wb.Worksheets.Item("List1")
$ws1.Activate()
$ws1.Range("A15").Value = 45494.636805555558
$ws1.Range("A15").Select()

# --- Sheet "Měření aktivity": fill in the row 16 measurement data ---
$ws2 = $wb.Worksheets.Item("Měření aktivity")
$ws2.Activate()

$ws2.Range("A16").Value = 45494.636805555558
$ws2.Range("B16").Value = 0.003
$ws2.Range("C16").Value = 0.003
$ws2.Range("D16").Value = 0.004
$ws2.Range("E16").Value = 0.004
$ws2.Range("F16").Value = 0.004
$ws2.Range("G16").Value = 226.5
$ws2.Range("H16").Value = 226.5
$ws2.Range("I16").Value = 226.6
$ws2.Range("J16").Value = 226.6
$ws2.Range("K16").Value = 226.6
$ws2.Range("L16").Value = 226.6
$ws2.Range("M16").Value = 226.6
$ws2.Range("N16").Value = 226.6
$ws2.Range("O16").Value = 226.6
$ws2.Range("P16").Value = 226.6

# Q16 already holds the formula
# =AVERAGE(G16:P16) - AVERAGE(Tabulka1[[#This Row],[č.1]:[č.5]])
# and will now evaluate instead of showing #DIV/0!

$ws2.Range("G17").Select()

$wb.Save()
